$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 7
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = 3
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = -4
$ws.Range("F27").Value = -2
$ws.Range("F28").Value = -2
$ws.Range("F30").Value = -2
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = -3
$ws.Range("F35").Value = 8
